$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update BOM "Value" column (column B) entries to reflect corrected / updated part values.
# Row 21 -> Reference "F2": value was the numeric "1206" (footprint code), correct it to "0.5A"
$ws.Range("B21").Value = "0.5A"

# Row 26 -> Reference "J9": value "S48B-ZROK-2A-R" renamed to "Yamaha 48p"
$ws.Range("B26").Value = "Yamaha 48p"

# Row 33 -> Reference "Q1,Q2,Q3,Q4,Q6,Q10,Q11,Q12": value "Q_NMOS_GDS" renamed to "BSP78/NCV8405A/NCV8406A"
$ws.Range("B33").Value = "BSP78/NCV8405A/NCV8406A"

# Row 47 -> Reference "R52,R53,R54,R55,R105,R106": value "15R" updated to "15R 2W"
$ws.Range("B47").Value = "15R 2W"

# Column B needs to widen (bestFit) to accommodate the new longer text values.
$ws.Columns("B:B").AutoFit() | Out-Null

# Scroll the sheet view so that row 41 is the top-left visible cell, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 41
